$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Step 1: shared-text header updates ----
$ws.Range("A8").Value = "Volume 30   Number  28"
$ws.Range("C9").Value = "Report Covering the Week  7/10/2023  Through  7/16/2023"

# ---- Step 2: type-switch cells (copy format+value from an untouched donor cell) ----
# text -> number (1): donor G14
$ws.Range("G14").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("G14").Copy()
$ws.Range("C14").PasteSpecial(-4163)
$ws.Range("G14").Copy()
$ws.Range("F14").PasteSpecial(-4122)
$ws.Range("G14").Copy()
$ws.Range("F14").PasteSpecial(-4163)
$ws.Range("G14").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("G14").Copy()
$ws.Range("C15").PasteSpecial(-4163)

# number -> text '0' (shared string 20): donor D14
$ws.Range("D14").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("D14").Copy()
$ws.Range("C17").PasteSpecial(-4163)
$ws.Range("D14").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("D14").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("D14").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("D14").Copy()
$ws.Range("C23").PasteSpecial(-4163)
$ws.Range("D14").Copy()
$ws.Range("F30").PasteSpecial(-4122)
$ws.Range("D14").Copy()
$ws.Range("F30").PasteSpecial(-4163)

# number -> text '***.*' (shared string 21): donor E14
$ws.Range("E14").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E22").PasteSpecial(-4163)

# text -> number (2): donor D17 (stays 2 throughout)
$ws.Range("D17").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("D17").Copy()
$ws.Range("D27").PasteSpecial(-4163)

# text -> number (0): donor E17 (captured BEFORE E17's own value edit below)
$ws.Range("E17").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("E17").Copy()
$ws.Range("E27").PasteSpecial(-4163)

# ---- Step 3: plain numeric value updates (style/type unchanged) ----
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 3
$ws.Range("K14").Value = 200
$ws.Range("L14").Value = 200
$ws.Range("M14").Value = 200
$ws.Range("N14").Value = 0
$ws.Range("F15").Value = 2
$ws.Range("I15").Value = 7
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 16.666666666666
$ws.Range("N15").Value = -58.823529411764
$ws.Range("C16").Value = 3
$ws.Range("E16").Value = -25
$ws.Range("F16").Value = 20
$ws.Range("G16").Value = 21
$ws.Range("H16").Value = -4.761904761904
$ws.Range("I16").Value = 106
$ws.Range("J16").Value = 128
$ws.Range("K16").Value = -17.1875
$ws.Range("L16").Value = 27.710843373494
$ws.Range("M16").Value = 60.606060606060
$ws.Range("N16").Value = -85.499316005472
$ws.Range("E17").Value = -100
$ws.Range("F17").Value = 14
$ws.Range("G17").Value = 14
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 105
$ws.Range("J17").Value = 93
$ws.Range("K17").Value = 12.903225806451
$ws.Range("L17").Value = 45.833333333333
$ws.Range("M17").Value = 123.404255319149
$ws.Range("N17").Value = -33.544303797468
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -80
$ws.Range("F18").Value = 11
$ws.Range("H18").Value = -52.173913043478
$ws.Range("I18").Value = 139
$ws.Range("J18").Value = 131
$ws.Range("K18").Value = 6.106870229007
$ws.Range("L18").Value = 20.869565217391
$ws.Range("M18").Value = 2.205882352941
$ws.Range("N18").Value = -90.014367816091
$ws.Range("C19").Value = 30
$ws.Range("D19").Value = 32
$ws.Range("E19").Value = -6.25
$ws.Range("F19").Value = 125
$ws.Range("G19").Value = 150
$ws.Range("H19").Value = -16.666666666666
$ws.Range("I19").Value = 879
$ws.Range("J19").Value = 898
$ws.Range("K19").Value = -2.115812917594
$ws.Range("L19").Value = 61.878453038674
$ws.Range("M19").Value = 31.194029850746
$ws.Range("N19").Value = -55.695564516129
$ws.Range("C20").Value = 6
$ws.Range("D20").Value = 11
$ws.Range("E20").Value = -45.454545454545
$ws.Range("F20").Value = 23
$ws.Range("G20").Value = 24
$ws.Range("H20").Value = -4.166666666666
$ws.Range("I20").Value = 98
$ws.Range("J20").Value = 100
$ws.Range("K20").Value = -2
$ws.Range("L20").Value = 25.641025641025
$ws.Range("M20").Value = 127.906976744186
$ws.Range("N20").Value = -94.624245748765
$ws.Range("C21").Value = 42
$ws.Range("D21").Value = 54
$ws.Range("E21").Value = -22.222222222222
$ws.Range("F21").Value = 196
$ws.Range("H21").Value = -15.879828326180
$ws.Range("I21").Value = 1337
$ws.Range("J21").Value = 1358
$ws.Range("K21").Value = -1.546391752577
$ws.Range("L21").Value = 48.720800889877
$ws.Range("M21").Value = 37.977296181630
$ws.Range("N21").Value = -78.110674525212
$ws.Range("F22").Value = 4
$ws.Range("H22").Value = 33.333333333333
$ws.Range("M22").Value = 100
$ws.Range("F23").Value = 1
$ws.Range("H23").Value = -66.666666666666
$ws.Range("C24").Value = 62
$ws.Range("D24").Value = 78
$ws.Range("E24").Value = -20.512820512820
$ws.Range("F24").Value = 286
$ws.Range("G24").Value = 351
$ws.Range("H24").Value = -18.518518518518
$ws.Range("I24").Value = 1747
$ws.Range("J24").Value = 2089
$ws.Range("K24").Value = -16.371469602680
$ws.Range("L24").Value = 32.851711026616
$ws.Range("M24").Value = 103.613053613054
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = -33.333333333333
$ws.Range("F25").Value = 22
$ws.Range("G25").Value = 29
$ws.Range("H25").Value = -24.137931034482
$ws.Range("I25").Value = 180
$ws.Range("J25").Value = 206
$ws.Range("K25").Value = -12.621359223301
$ws.Range("L25").Value = 7.142857142857
$ws.Range("M25").Value = -9.090909090909
$ws.Range("I26").Value = 16
$ws.Range("K26").Value = 60
$ws.Range("L26").Value = 60
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 25
$ws.Range("I27").Value = 54
$ws.Range("J27").Value = 46
$ws.Range("K27").Value = 17.391304347826
$ws.Range("L27").Value = 17.391304347826
$ws.Range("G30").Value = 3
$ws.Range("H30").Value = -100
